$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are plain text even when they look numeric
# (e.g. "41.539.77", "0.180"). Setting .Value directly on a numeric-looking
# string makes Excel coerce it to a real number (losing formatting like a
# trailing zero), so those cells are briefly switched to Text format, the
# text is written, and then the original (unstyled) look is restored by
# copying the style of the untouched Link cell on the same row.

$ws.Range("D2").Value = '41.539.77'
$ws.Range("E2").Value = '  -2.50%  '
$ws.Range("D3").Value = '2.463.54'
$ws.Range("E3").Value = '  -2.72%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.39'
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.88'
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = '  -5.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = '  -2.97%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.81'
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = '  -5.95%  '
$ws.Range("E11").Value = '  -2.73%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.99'
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = '  -4.80%  '
$ws.Range("D14").Value = '2.846.15'
$ws.Range("E14").Value = '  -3.28%  '
$ws.Range("D15").Value = '2.466.04'
$ws.Range("E15").Value = '  -2.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.59'
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = '  -8.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.789'
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = '  -3.82%  '
$ws.Range("D18").Value = '41.537.16'
$ws.Range("E18").Value = '  -2.45%  '
$ws.Range("E19").Value = '  -6.59%  '
$ws.Range("D20").Value = '0.0₃0917'
$ws.Range("E20").Value = '  -3.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.61'
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.48'
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.99'
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = '  -2.48%  '
$ws.Range("E24").Value = '  -4.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.95'
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = '  -4.76%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -4.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = '  -4.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.75'
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = '  -3.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.36'
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = '  -7.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.20'
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.65'
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.63'
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -8.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0756'
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = '  -4.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.02'
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = '  -4.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.16'
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = '  -6.24%  '
$ws.Range("E38").Value = '  -7.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.105'
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = '  -5.41%  '
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.04'
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  -6.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.51'
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").Value = '1.987.19'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("E45").Value = '  -4.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.06'
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = '  -7.65%  '
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.66'
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = '  -5.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '69.92'
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = '  -3.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.69'
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.180'
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = '  -6.05%  '
